$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.065.75"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "1.668.35"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "216.17"
$ws.Range("E5").Value = "  -1.16%  "

$ws.Range("D6").Value = "0.5116"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "0.2688"
$ws.Range("E8").Value = "  +0.93%  "

$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").Value = "21.79"
$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("D11").Value = "0.07441"
$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").Value = "1.705.30"
$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").Value = "4.510"
$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("D14").Value = "0.5812"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").Value = "0.000008487"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").Value = "64.07"
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "25.872.34"
$ws.Range("E17").Value = "  -1.67%  "

$ws.Range("D18").Value = "4.927"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").Value = "10.79"
$ws.Range("E20").Value = "  -0.66%  "

$ws.Range("D21").Value = "189.24"
$ws.Range("E21").Value = "  +1.37%  "

$ws.Range("D22").Value = "6.182"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "144.51"
$ws.Range("E24").Value = "  +0.82%  "

$ws.Range("D25").Value = "7.593"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").Value = "0.1226"
$ws.Range("E26").Value = "  +4.02%  "

$ws.Range("D27").Value = "15.68"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").Value = "0.06656"
$ws.Range("E28").Value = "  +14.70%  "

$ws.Range("D29").Value = "1.340"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("D30").Value = "1.313"
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").Value = "3.572"
$ws.Range("E31").Value = "  +1.71%  "

$ws.Range("D32").Value = "3.524"
$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("D33").Value = "1.661"
$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("D34").Value = "1.016"
$ws.Range("E34").Value = "  +1.19%  "

$ws.Range("D35").Value = "0.6151"
$ws.Range("E35").Value = "  +3.36%  "

$ws.Range("D36").Value = "2.367"
$ws.Range("E36").Value = "  +0.24%  "

$ws.Range("D37").Value = "2.688"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("E38").Value = "  +6.10%  "

$ws.Range("D39").Value = "1.093.59"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").Value = "0.01596"
$ws.Range("E40").Value = "  -0.71%  "

$ws.Range("D41").Value = "0.8700"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").Value = "101.07"
$ws.Range("E43").Value = "  +1.41%  "

$ws.Range("D44").Value = "1.815.30"
$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("E45").Value = "  +1.79%  "

$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.127"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("D49").Value = "0.05236"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("D50").Value = "0.4280"
$ws.Range("E50").Value = "  -0.81%  "

$ws.Range("D51").Value = "5.981"
$ws.Range("E51").Value = "  +2.34%  "
